# Refresh the cryptos price table (columns D = Price, E = Volume(1h))
# with the latest values, cell by cell, mirroring the upstream scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.402.76'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.570.91'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '''1.002'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").Value = '''289.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '''0.3742'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '''49.40'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").Value = '''0.3363'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("E10").Value = '  -3.14%  '
$ws.Range("D11").Value = '''0.07416'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.95%  '
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '''20.89'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.35%  '
$ws.Range("D14").Value = '''5.899'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = '''6.854'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.22%  '
$ws.Range("D16").Value = '1.574.67'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").Value = '''0.00001108'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.95%  '
$ws.Range("D18").Value = '''89.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.08%  '
$ws.Range("D19").Value = '''0.06668'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("D21").Value = '''6.161'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.22%  '
$ws.Range("D22").Value = '''16.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.41%  '
$ws.Range("D23").Value = '''11.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").Value = '22.392.34'
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").Value = '''2.362'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("D26").Value = '''2.536'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.85%  '
$ws.Range("D27").Value = '''19.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("D28").Value = '''147.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("D29").Value = '''4.993'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("D30").Value = '''124.73'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("D31").Value = '1.745.69'
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").Value = '''0.9976'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.20%  '
$ws.Range("E33").Value = '  -1.38%  '
$ws.Range("D34").Value = '''5.933'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.34%  '
$ws.Range("D35").Value = '''9.664'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.90%  '
$ws.Range("D36").Value = '''0.08380'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.37%  '
$ws.Range("D37").Value = '''1.372'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.63%  '
$ws.Range("D38").Value = '''0.02455'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.88%  '
$ws.Range("D39").Value = '''0.2244'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.34%  '
$ws.Range("D40").Value = '''0.06388'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").Value = '''5.362'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.97%  '
$ws.Range("D42").Value = '''11.09'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.38%  '
$ws.Range("D43").Value = '''0.6183'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.92%  '
$ws.Range("D44").Value = '''1.003'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.28%  '
$ws.Range("D45").Value = '''13.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").Value = '''3.801'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("D47").Value = '''0.5777'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("D48").Value = '''2.048'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("D49").Value = '''125.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("E50").Value = '  -3.45%  '
$ws.Range("D51").Value = '''0.07297'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.39%  '
